$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-6 ---
$ws.Range("C2").Value = 5
$ws.Range("G2").Value = 0.5028938270292095
$ws.Range("H2").Value = 0.00407062
$ws.Range("I2").Value = 0.6900002850299007
$ws.Range("J2").Value = 0.00210366

$ws.Range("C3").Value = 5

$ws.Range("C4").Value = 5
$ws.Range("G4").Value = 0.4846250348041683
$ws.Range("H4").Value = 0.00568506
$ws.Range("I4").Value = 0.7150502036319386
$ws.Range("J4").Value = 0.00392709

$ws.Range("B5").Value = "DT_sample:1356_cf:1_mean_of_5_iterations"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 1356
$ws.Range("F5").Value = 33003
$ws.Range("G5").Value = 0.4672272462712558
$ws.Range("H5").Value = 0.008394789999999999
$ws.Range("I5").Value = 0.7008665457116703
$ws.Range("J5").Value = 0.00599623

$ws.Range("B6").Value = "DT_sample:1356_cf:2_mean_of_5_iterations"
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 1356
$ws.Range("F6").Value = 34359
$ws.Range("G6").Value = 0.4610902867966988
$ws.Range("H6").Value = 0.00714374
$ws.Range("I6").Value = 0.6934237462795376
$ws.Range("J6").Value = 0.00477817

# --- Add new rows 7-10, copying style of column A from row 2 ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "DT_sample:2712_cf:1_mean_of_5_iterations"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 2712
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 34359
$ws.Range("G7").Value = 0.4637420381395428
$ws.Range("H7").Value = 0.00327374
$ws.Range("I7").Value = 0.6960786134391599
$ws.Range("J7").Value = 0.00174809

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "DT_sample:2712_cf:2_mean_of_5_iterations"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 2712
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 37071
$ws.Range("G8").Value = 0.4586469228887627
$ws.Range("H8").Value = 0.00432843
$ws.Range("I8").Value = 0.6899385212429496
$ws.Range("J8").Value = 0.00272537

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "DT_sample:4069_cf:1_mean_of_5_iterations"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 4069
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 35716
$ws.Range("G9").Value = 0.4655479305685267
$ws.Range("H9").Value = 0.01236409
$ws.Range("I9").Value = 0.6949310611334323
$ws.Range("J9").Value = 0.0067427

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "DT_sample:4069_cf:2_mean_of_5_iterations"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 4069
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 39785
$ws.Range("G10").Value = 0.4577977849919957
$ws.Range("H10").Value = 0.00696073
$ws.Range("I10").Value = 0.6864653332383234
$ws.Range("J10").Value = 0.00334761

$excel.CutCopyMode = $false
Write-Host "edit applied"